$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 24500
$ws.Range("J21").Value = 24500
$ws.Range("L21").Value = 24500
$ws.Range("N21").Value = -25436

$ws.Range("H23").Value = 24500
$ws.Range("J23").Value = 24500
$ws.Range("L23").Value = 24500
$ws.Range("N23").Value = -24968

$ws.Range("H33").Value = 436.1875
$ws.Range("I33").Value = 244.125
$ws.Range("J33").Value = 1012.375
$ws.Range("K33").Value = 244.125
$ws.Range("L33").Value = 1012.375
$ws.Range("M33").Value = -15.125
$ws.Range("N33").Value = -1470.375

$ws.Range("H40").Value = 3035.9062
$ws.Range("J40").Value = 3512.2173
$ws.Range("L40").Value = 3512.2173
$ws.Range("N40").Value = -3862.2173

$ws.Range("H43").Value = 1800
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 1900
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 1900
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -2038

$ws.Range("H51").Value = 3099.5
$ws.Range("J51").Value = 2998.5715
$ws.Range("L51").Value = 2998.5715
$ws.Range("N51").Value = -3966.5715

$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 500
$ws.Range("M94").Value = -49

$ws.Range("H101").Value = 934.125
$ws.Range("I101").Value = 1128
$ws.Range("J101").Value = 611
$ws.Range("K101").Value = 3384
$ws.Range("L101").Value = 1833
$ws.Range("M101").Value = -1762
$ws.Range("N101").Value = -5077

$ws.Range("H111").Value = 1249
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 1998
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 5994
$ws.Range("M111").Value = 1567
$ws.Range("N111").Value = -12128

$ws.Range("H132").Value = 8139.84
$ws.Range("I132").Value = 9536.895
$ws.Range("J132").Value = 3715.8333
$ws.Range("K132").Value = 28610.685
$ws.Range("L132").Value = 11147.4999
$ws.Range("M132").Value = -26080.685
$ws.Range("N132").Value = -16207.4999

$ws.Range("H138").Value = 3537.3967
$ws.Range("J138").Value = 4849.0347
$ws.Range("L138").Value = 14547.1041
$ws.Range("N138").Value = -24827.1041

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2715.5386
$ws.Range("I2").Value = 2221.8572
$ws.Range("K2").Value = 2221.8572
$ws.Range("M2").Value = -2108.8572

$ws.Range("H32").Value = 28128608
$ws.Range("I32").Value = 15386748
$ws.Range("J32").Value = 83343336
$ws.Range("K32").Value = 15386748
$ws.Range("L32").Value = 83343336
$ws.Range("M32").Value = -15386461
$ws.Range("N32").Value = -83343910

$ws.Range("H45").Value = 3070.7827
$ws.Range("I45").Value = 2774.3333
$ws.Range("K45").Value = 2774.3333
$ws.Range("M45").Value = -2397.3333

$ws.Range("H55").Value = 50021.5

$ws.Range("H102").Value = 5333
$ws.Range("I102").Value = 5999.5
$ws.Range("K102").Value = 5999.5
$ws.Range("M102").Value = -4377.5

$ws.Range("H116").Value = 2715.5386
$ws.Range("I116").Value = 2221.8572
$ws.Range("K116").Value = 2221.8572
$ws.Range("M116").Value = 72.14280000000008

$ws.Range("H122").Value = 6057.645
$ws.Range("I122").Value = 5076.6816
$ws.Range("K122").Value = 15230.0448
$ws.Range("M122").Value = -12780.0448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2715.5386
$ws.Range("I3").Value = 2221.8572
$ws.Range("K3").Value = 2221.8572
$ws.Range("M3").Value = -2107.8572

$ws.Range("H22").Value = 278.85715
$ws.Range("I22").Value = 278.85715
$ws.Range("K22").Value = 278.85715
$ws.Range("M22").Value = -105.85715

$ws.Range("H86").Value = 4680.3716
$ws.Range("I86").Value = 4015.7693
$ws.Range("J86").Value = 6600.3335
$ws.Range("K86").Value = 4015.7693
$ws.Range("L86").Value = 6600.3335
$ws.Range("M86").Value = -2892.7693
$ws.Range("N86").Value = -8846.333500000001

$ws.Range("H89").Value = 4680.3716
$ws.Range("I89").Value = 4015.7693
$ws.Range("J89").Value = 6600.3335
$ws.Range("K89").Value = 20078.8465
$ws.Range("L89").Value = 33001.6675
$ws.Range("M89").Value = -14462.8465
$ws.Range("N89").Value = -44233.6675

$ws.Range("H99").Value = 3841.2856
$ws.Range("I99").Value = 3777.8
$ws.Range("K99").Value = 3777.8
$ws.Range("M99").Value = -2279.8

$ws.Range("H105").Value = 2198.7646
$ws.Range("I105").Value = 1958.6
$ws.Range("K105").Value = 1958.6
$ws.Range("M105").Value = -211.5999999999999

$ws.Range("H107").Value = 1050.0769
$ws.Range("I107").Value = 867.5
$ws.Range("J107").Value = 1658.6666
$ws.Range("K107").Value = 867.5
$ws.Range("L107").Value = 1658.6666
$ws.Range("M107").Value = 1052.5
$ws.Range("N107").Value = -5498.6666

$ws.Range("H132").Value = 112140.766
$ws.Range("J132").Value = 112140.766
$ws.Range("L132").Value = 112140.766
$ws.Range("N132").Value = -122260.766

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4327.2905
$ws.Range("I31").Value = 3194.9473
$ws.Range("K31").Value = 3194.9473
$ws.Range("M31").Value = -2899.9473

$ws.Range("H34").Value = 4327.2905
$ws.Range("I34").Value = 3194.9473
$ws.Range("K34").Value = 3194.9473
$ws.Range("M34").Value = -2992.9473

$ws.Range("H75").Value = 113984.5
$ws.Range("J75").Value = 113984.5
$ws.Range("L75").Value = 113984.5
$ws.Range("N75").Value = -115980.5

$ws.Range("H78").Value = 113984.5
$ws.Range("J78").Value = 113984.5
$ws.Range("L78").Value = 341953.5
$ws.Range("N78").Value = -351937.5

$ws.Range("H86").Value = 3202.8125
$ws.Range("I86").Value = 3270
$ws.Range("K86").Value = 3270
$ws.Range("M86").Value = -2147

$ws.Range("H89").Value = 3202.8125
$ws.Range("I89").Value = 3270
$ws.Range("K89").Value = 16350
$ws.Range("M89").Value = -10734

$ws.Range("H99").Value = 1953.7142
$ws.Range("I99").Value = 1996.75
$ws.Range("K99").Value = 1996.75
$ws.Range("M99").Value = -498.75

$ws.Range("H105").Value = 1613.7222
$ws.Range("I105").Value = 980.8461
$ws.Range("K105").Value = 980.8461
$ws.Range("M105").Value = 766.1539

$ws.Range("H124").Value = 33269.445
$ws.Range("J124").Value = 32178.125
$ws.Range("L124").Value = 32178.125
$ws.Range("N124").Value = -37088.125

$ws.Range("H126").Value = 1953.7142
$ws.Range("I126").Value = 1996.75
$ws.Range("K126").Value = 5990.25
$ws.Range("M126").Value = -3520.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8494.5625
$ws.Range("I3").Value = 8193
$ws.Range("J3").Value = 8997.166999999999
$ws.Range("K3").Value = 24579
$ws.Range("L3").Value = 26991.501
$ws.Range("M3").Value = -24467
$ws.Range("N3").Value = -27215.501

$ws.Range("H4").Value = 37982110
$ws.Range("I4").Value = 65227364
$ws.Range("J4").Value = 24569064
$ws.Range("K4").Value = 195682092
$ws.Range("L4").Value = 73707192
$ws.Range("M4").Value = -195681980
$ws.Range("N4").Value = -73707416

$ws.Range("H44").Value = 791.5
$ws.Range("I44").Value = 937.25
$ws.Range("J44").Value = 500
$ws.Range("K44").Value = 2811.75
$ws.Range("L44").Value = 1500
$ws.Range("M44").Value = -2413.75
$ws.Range("N44").Value = -2296

$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H138").Value = 26627286
$ws.Range("J138").Value = 53252748
$ws.Range("L138").Value = 159758244
$ws.Range("N138").Value = -159768524

$ws.Range("H140").Value = 15386083
$ws.Range("I140").Value = 22223370
$ws.Range("J140").Value = 2187.5
$ws.Range("K140").Value = 66670110
$ws.Range("L140").Value = 6562.5
$ws.Range("M140").Value = -66664930
$ws.Range("N140").Value = -16922.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 20021
$ws.Range("I38").Value = 20021
$ws.Range("K38").Value = 20021
$ws.Range("M38").Value = -19558

$ws.Range("H102").Value = 1932.48
$ws.Range("J102").Value = 1946
$ws.Range("L102").Value = 1946
$ws.Range("N102").Value = -5190

$ws.Range("H122").Value = 4840.7646
$ws.Range("I122").Value = 5018.625
$ws.Range("K122").Value = 15055.875
$ws.Range("M122").Value = -12605.875

$ws.Range("H126").Value = 1446.9
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230

$ws.Range("H128").Value = 152950
$ws.Range("J128").Value = 152950
$ws.Range("L128").Value = 152950
$ws.Range("N128").Value = -162910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 96141.28999999999
$ws.Range("J128").Value = 96141.28999999999
$ws.Range("L128").Value = 96141.28999999999
$ws.Range("N128").Value = -106101.29

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12108.333
$ws.Range("J41").Value = 12108.333
$ws.Range("L41").Value = 12108.333
$ws.Range("N41").Value = -12888.333
